$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*reside in the County of*") {
        $target = $p
        break
    }
}

$r = $target.Range.Duplicate
$found = $r.Find.Execute("{{COUNTY}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$prevSpace = $d.Range($r.Start - 1, $r.Start)
$prevSpace.Text = ""

$r = $target.Range.Duplicate
$found = $r.Find.Execute("{{COUNTY}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = " Frio"

$spaceRange = $d.Range($r.Start, $r.Start + 1)
$wordRange  = $d.Range($r.Start + 1, $r.End)

$spaceRange.Bold = 0
$wordRange.Bold = 0

Write-Output "after Bold=0 on both:"
